$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only columns B:F (the error-metric columns) shift down by one row for rows 2-11;
# column A (quarter labels) and column G (N) stay fixed per row.
# Shift existing values in B2:F10 down into B3:F11 (working bottom-up to avoid overwrite).
for ($r = 10; $r -ge 2; $r--) {
    $src = $ws.Range("B" + $r + ":F" + $r)
    $dst = $ws.Range("B" + ($r + 1) + ":F" + ($r + 1))
    $dst.Value2 = $src.Value2
}

# Populate the new row 2 (B2:F2) with the new values from the diff.
$ws.Range("B2").Value2 = 0.1560865643779764
$ws.Range("C2").Value2 = 0.3515450347245845
$ws.Range("D2").Value2 = 0.2150214299408537
$ws.Range("E2").Value2 = 0.4637040326985023
$ws.Range("F2").Value2 = 0.446457955381491
